# Finish first complete draft and spell check:
# adjust heading spacing on the built-in Heading 1/2/3 styles.
$d = $word.ActiveDocument

$h1 = $d.Styles("Heading1")
$h1.ParagraphFormat.SpaceBefore = 18
$h1.ParagraphFormat.SpaceAfter = 12

$h2 = $d.Styles("Heading2")
$h2.ParagraphFormat.SpaceAfter = 6

$h3 = $d.Styles("Heading3")
$h3.ParagraphFormat.SpaceAfter = 0
